$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = "9665668010"
$ws.Range("T2").Value = "16/08/2016"
$ws.Range("AX2").Value = "9665668010"

$ws.Range("S3").Value = "9665668010"
$ws.Range("T3").Value = "18/08/2016"
$ws.Range("AX3").Value = "9665668010"

$ws.Range("S4").Value = "9665668010"
$ws.Range("T4").Value = "31/08/2016"
$ws.Range("AX4").Value = "9665668010"
